$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Tuesday's "Regular Hrs" entry from 1 to 2 (dependent totals/formulas
# recalc automatically: daily total, weekly total, pay, and summary cells).
$ws.Range("B12").Value = 2

# Move the active selection to the cell that was edited.
$ws.Range("B12").Select() | Out-Null
